$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.175.24"
$ws.Range("E2").Value = "  +3.08%  "

$ws.Range("D3").Value = "1.580.26"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("E6").Value = "  +5.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.89%  "

$ws.Range("E9").Value = "  +2.36%  "

$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").Value = "1.806.11"
$ws.Range("E12").Value = "  +1.80%  "

$ws.Range("D13").Value = "1.579.84"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").Value = "29.210.71"
$ws.Range("E14").Value = "  +3.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.44%  "

$ws.Range("E16").Value = "  +2.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("E19").Value = "  +1.48%  "

$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.108"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.50%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +1.54%  "

$ws.Range("D33").Value = "1.423.20"
$ws.Range("E33").Value = "  +2.53%  "

$ws.Range("E34").Value = "  +1.53%  "

$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.16%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("E40").Value = "  +3.57%  "

$ws.Range("E41").Value = "  +2.52%  "

$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "53.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.01%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("E44").Value = "  +1.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "

$ws.Range("D48").Value = "1.718.41"
$ws.Range("E48").Value = "  +1.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.842"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  +0.35%  "
